$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.858.48"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "3.157.89"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'531.89"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'140.34"
$ws.Range("E6").Value = "  +1.49%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "'0.538"
$ws.Range("E8").Value = "  +16.59%  "
$ws.Range("D9").Value = "'7.35"
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").Value = "'0.436"
$ws.Range("E10").Value = "  +6.74%  "
$ws.Range("E11").Value = "  +3.40%  "
$ws.Range("E12").Value = "  +2.83%  "
$ws.Range("D13").Value = "3.697.90"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").Value = "'26.26"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("E15").Value = "  +6.27%  "
$ws.Range("D16").Value = "58.886.02"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "'6.28"
$ws.Range("E17").Value = "  +4.93%  "
$ws.Range("D18").Value = "3.156.65"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").Value = "'13.09"
$ws.Range("E19").Value = "  +3.70%  "
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").Value = "'374.12"
$ws.Range("E21").Value = "  +5.75%  "
$ws.Range("D22").Value = "'5.81"
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'70.30"
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("D25").Value = "'0.522"
$ws.Range("E25").Value = "  +3.84%  "
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "'8.28"
$ws.Range("E28").Value = "  +13.91%  "
$ws.Range("D29").Value = "0.0₃0867"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'22.29"
$ws.Range("E30").Value = "  +4.68%  "
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").Value = "'6.15"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.89"
$ws.Range("E32").Value = "  +1.22%  "
$ws.Range("D33").Value = "'5.21"
$ws.Range("E33").Value = "  +4.16%  "
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "'6.30"
$ws.Range("E35").Value = "  +3.98%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'158.94"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").Value = "'1.34"
$ws.Range("E37").Value = "  +7.02%  "
$ws.Range("D38").Value = "'25.24"
$ws.Range("E38").Value = "  -2.80%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.650.26"
$ws.Range("E40").Value = "  +10.46%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "'0.0687"
$ws.Range("E41").Value = "  +2.41%  "
$ws.Range("D42").Value = "'4.26"
$ws.Range("E42").Value = "  +6.58%  "
$ws.Range("D43").Value = "'0.0289"
$ws.Range("E43").Value = "  +7.65%  "
$ws.Range("D44").Value = "'38.85"
$ws.Range("E44").Value = "  +3.32%  "
$ws.Range("D45").Value = "'0.712"
$ws.Range("E45").Value = "  +1.87%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "3.196.49"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("D48").Value = "'0.103"
$ws.Range("E48").Value = "  +14.08%  "
$ws.Range("D49").Value = "'0.988"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").Value = "'6.21"
$ws.Range("E50").Value = "  +3.07%  "
$ws.Range("D51").Value = "'20.29"
$ws.Range("E51").Value = "  +2.35%  "
